$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pages")

# --- Row 13 ("User list"): mark as Done in the Notes column (H13) ---
$ws.Range("H13").Value = "Done"

# --- Row 14 ("User RUD"): flag row red and note it is "Not persent" ---
$redCells14 = @("A14", "B14", "D14", "E14", "G14", "H14")
foreach ($addr in $redCells14) {
    $ws.Range($addr).Interior.Color = 255
}
$ws.Range("H14").Value = "Not persent"

# --- Row 16 ("Store support requests List"): flag row red and note "Not present" ---
$redCells16 = @("A16", "B16", "D16", "E16", "G16", "H16")
foreach ($addr in $redCells16) {
    $ws.Range($addr).Interior.Color = 255
}
$ws.Range("H16").Value = "Not present"

# --- Row 18 ("Store Logs"): flag row red and note "Not present" ---
$redCells18 = @("A18", "B18", "D18", "E18", "G18", "H18")
foreach ($addr in $redCells18) {
    $ws.Range($addr).Interior.Color = 255
}
$ws.Range("H18").Value = "Not present"

# --- Update the saved selection on the "Pages" sheet to F9 (was F11) ---
$ws.Range("F9").Select()

$wb.Save()
